$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark several existing projects as "Added to map?" = TRUE before removing the
# HS2 Demonstration London row (row 41).
$ws.Range("C6").Value = $true
$ws.Range("C37").Value = $true
$ws.Range("C64").Value = $true
$ws.Range("C82").Value = $true
$ws.Range("C83").Value = $true
$ws.Range("C84").Value = $true

# Remove the "HS2 Demonstration London" entry entirely (row 41), shifting
# everything below it up by one row.
$ws.Rows.Item(41).Delete()

# Restore a sensible scroll position/selection similar to the saved file.
[void]$ws.Application.Goto($ws.Range("D73"), $false)
[void]$ws.Range("D73").Select()
